# Refresh the live crypto price/volume snapshot in Sheet1, as produced by
# the scheduled "Updated cryptos list ... with GitHub Actions" job.
#
# Most updated cells are plain text (Price/Volume columns are formatted
# strings, not numbers) so a straightforward Range.Value assignment is
# enough. A handful of new Price values are plain-looking decimals (e.g.
# "603.89"); Excel would otherwise auto-convert those to numeric cells on
# assignment (and mangle the decimal via float rounding). For those we
# write the value with a leading apostrophe (forces text, like a user
# typing it in the UI) and then reapply the "Normal" style so the cell
# format index is left exactly as it was (no lingering quote-prefix style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.640.04'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '3.297.68'
$ws.Range('E3').Value = '  +5.56%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''603.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').Value = '''141.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.26%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.296.70'
$ws.Range('E8').Value = '  +5.78%  '
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('E11').Value = '  +4.57%  '
$ws.Range('D12').Value = '''0.471'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').Value = '''34.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').Value = '3.839.57'
$ws.Range('E15').Value = '  +5.65%  '
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '3.292.94'
$ws.Range('E17').Value = '  +5.30%  '
$ws.Range('D18').Value = '63.707.24'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '''6.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('D20').Value = '''479.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Value = '''14.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '''0.732'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.26%  '
$ws.Range('D23').Value = '''8.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.77%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''13.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.57%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''84.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').Value = '''7.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.27%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  +3.46%  '
$ws.Range('E31').Value = '  +3.98%  '
$ws.Range('D32').Value = '''29.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.09%  '
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('E34').Value = '  +0.40%  '
$ws.Range('D36').Value = '''5.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.34%  '
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('D38').Value = '0.0₃0745'
$ws.Range('E38').Value = '  +7.39%  '
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('D40').Value = '''425.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').Value = '3.052.03'
$ws.Range('E41').Value = '  +5.32%  '
$ws.Range('D42').Value = '''8.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').Value = '''0.266'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('D46').Value = '''2.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.38%  '
$ws.Range('D47').Value = '''26.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.48%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('D50').Value = '''2.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('D51').Value = '''124.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.44%  '
